# CIERRE 19 MAR 22
# Add a new "SEMANA 11" week sheet (duplicated from the SEMANA 10 sheet),
# insert a fresh blank "Hoja1" sheet, and refresh the cash-breakdown sheet.

$wb = $excel.ActiveWorkbook

$week10 = $wb.Worksheets.Item("SEMANA   10    2022   ")
$hoja2  = $wb.Worksheets.Item("Hoja2")
$hoja5  = $wb.Worksheets.Item("Hoja5")
$hoja3  = $wb.Worksheets.Item("Hoja3")

# --- 1. Build the new "SEMANA 11" sheet as a duplicate of the SEMANA 10 sheet ---
$week10.Copy($null, $week10)
$week11 = $wb.Worksheets.Item("SEMANA   10    2022    (2)")

# Drop the now-unused placeholder sheet and rename the duplicate into place.
$hoja2.Delete() | Out-Null
$week11.Name = "SEMANA   11     2022   "

# Update the week banner text for SEMANA 11.
$week11.Range("B1").Value = "SEMANA #  11    DEL     14   AL  20   MARZO       2 0 2 2 "

# TEODORA ARELLANO PEREZ worked 3 extra hours this week; update her row and
# the manual cash-denomination breakdown used to pay her.
$week11.Range("F14").Value = 3
$week11.Range("T14").Value = 6
$week11.Range("W14").Value = 3
$week11.Range("X14").Value = 5

# View state: this is now the freshly-edited / focused sheet.
$week11.Range("P1").Select() | Out-Null
$week11.Range("W21").Select() | Out-Null

# --- 2. Remove the SEMANA 10 sheet's red tab color; it's no longer the latest week ---
$week10.Tab.ColorIndex = -4142
$week10.Range("R22").Select() | Out-Null

# --- 3. Insert a brand-new blank sheet "Hoja1" right after SEMANA 11 ---
$hoja1 = $wb.Worksheets.Add($null, $week11)
$hoja1.Name = "Hoja1"

# --- 4. Hoja5: selection moves to B4 (no longer the active tab) ---
$hoja5.Range("B4").Select() | Out-Null

# --- 5. Hoja3: refresh the denomination counts to match the new totals ---
$hoja3.Range("F6").Value = 18
$hoja3.Range("F9").Value = 3
$hoja3.Range("F10").Value = 8
$hoja3.Range("G19").Select() | Out-Null

# --- 6. Make SEMANA 11 the active tab, as in the committed workbook ---
$week11.Activate()
$week11.Range("W21").Select() | Out-Null
